# Auto-generated edit script.
# Applies the "Automatic update of files" change to the Artfynd (SSF)
# export sheet: rows 2-22 are re-ordered (the source report re-ran with a
# different sort/paging of the same underlying sightings) and the taxon
# "Taxonsorteringsordning" lookup values (column B) plus a couple of
# downstream summary rows (23-24) are refreshed to their new values.
# Concretely this rewrites each touched cell to the value it holds after
# the update, which reproduces the permutation + lookup refresh exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet



# Row 2
$ws.Cells.Item(2,1).Value = 111896690  # A2
$ws.Cells.Item(2,2).Value = 90835  # B2
$ws.Cells.Item(2,5).Value = 5964  # E2
$ws.Cells.Item(2,6).Value = 'Fjällig taggsvamp s.str.'  # F2
$ws.Cells.Item(2,7).Value = 'Sarcodon imbricatus s.str.'  # G2
$ws.Cells.Item(2,8).Value = '(L.:Fr.) P.Karst.'  # H2
$ws.Cells.Item(2,17).Value = 575060  # Q2
$ws.Cells.Item(2,18).Value = 6703377  # R2

# Row 3
$ws.Cells.Item(3,1).Value = 111896637  # A3
$ws.Cells.Item(3,2).Value = 90480  # B3
$ws.Cells.Item(3,17).Value = 575088  # Q3
$ws.Cells.Item(3,18).Value = 6703396  # R3

# Row 4
$ws.Cells.Item(4,1).Value = 111896652  # A4
$ws.Cells.Item(4,2).Value = 89331  # B4
$ws.Cells.Item(4,5).Value = 3215  # E4
$ws.Cells.Item(4,6).Value = 'Rödgul trumpetsvamp'  # F4
$ws.Cells.Item(4,7).Value = 'Craterellus lutescens'  # G4
$ws.Cells.Item(4,8).Value = '(Fr.) Fr.'  # H4
$ws.Cells.Item(4,17).Value = 575067  # Q4
$ws.Cells.Item(4,18).Value = 6703456  # R4

# Row 5
$ws.Cells.Item(5,1).Value = 111896644  # A5
$ws.Cells.Item(5,2).Value = 90480  # B5
$ws.Cells.Item(5,5).Value = 4769  # E5
$ws.Cells.Item(5,6).Value = 'Svavelriska'  # F5
$ws.Cells.Item(5,7).Value = 'Lactarius scrobiculatus'  # G5
$ws.Cells.Item(5,8).Value = '(Scop.:Fr.) Fr.'  # H5
$ws.Cells.Item(5,17).Value = 575036  # Q5
$ws.Cells.Item(5,18).Value = 6703432  # R5

# Row 6
$ws.Cells.Item(6,1).Value = 111883983  # A6
$ws.Cells.Item(6,2).Value = 90480  # B6
$ws.Cells.Item(6,4).Value = 'LC'  # D6
$ws.Cells.Item(6,5).Value = 4769  # E6
$ws.Cells.Item(6,6).Value = 'Svavelriska'  # F6
$ws.Cells.Item(6,7).Value = 'Lactarius scrobiculatus'  # G6
$ws.Cells.Item(6,8).Value = '(Scop.:Fr.) Fr.'  # H6
$ws.Cells.Item(6,17).Value = 575058  # Q6
$ws.Cells.Item(6,18).Value = 6703446  # R6

# Row 7
$ws.Cells.Item(7,1).Value = 111896642  # A7
$ws.Cells.Item(7,2).Value = 90480  # B7
$ws.Cells.Item(7,17).Value = 575014  # Q7
$ws.Cells.Item(7,18).Value = 6703387  # R7

# Row 8
$ws.Cells.Item(8,1).Value = 111896633  # A8
$ws.Cells.Item(8,2).Value = 90480  # B8
$ws.Cells.Item(8,5).Value = 4769  # E8
$ws.Cells.Item(8,6).Value = 'Svavelriska'  # F8
$ws.Cells.Item(8,7).Value = 'Lactarius scrobiculatus'  # G8
$ws.Cells.Item(8,8).Value = '(Scop.:Fr.) Fr.'  # H8
$ws.Cells.Item(8,16).Value = 'Kratte masugn, Gstr'  # P8
$ws.Cells.Item(8,17).Value = 575100  # Q8
$ws.Cells.Item(8,18).Value = 6703444  # R8
$ws.Cells.Item(8,49).Value = 'Philipp Weiss'  # AW8
$ws.Cells.Item(8,50).Value = 'Philipp Weiss'  # AX8

# Row 9
$ws.Cells.Item(9,1).Value = 111896634  # A9
$ws.Cells.Item(9,2).Value = 90480  # B9
$ws.Cells.Item(9,5).Value = 4769  # E9
$ws.Cells.Item(9,6).Value = 'Svavelriska'  # F9
$ws.Cells.Item(9,7).Value = 'Lactarius scrobiculatus'  # G9
$ws.Cells.Item(9,8).Value = '(Scop.:Fr.) Fr.'  # H9
$ws.Cells.Item(9,17).Value = 575048  # Q9
$ws.Cells.Item(9,18).Value = 6703452  # R9

# Row 10
$ws.Cells.Item(10,1).Value = 111896653  # A10
$ws.Cells.Item(10,2).Value = 89331  # B10
$ws.Cells.Item(10,5).Value = 3215  # E10
$ws.Cells.Item(10,6).Value = 'Rödgul trumpetsvamp'  # F10
$ws.Cells.Item(10,7).Value = 'Craterellus lutescens'  # G10
$ws.Cells.Item(10,8).Value = '(Fr.) Fr.'  # H10
$ws.Cells.Item(10,17).Value = 575075  # Q10
$ws.Cells.Item(10,18).Value = 6703404  # R10

# Row 11
$ws.Cells.Item(11,1).Value = 111884471  # A11
$ws.Cells.Item(11,2).Value = 89047  # B11
$ws.Cells.Item(11,4).Value = 'NT'  # D11
$ws.Cells.Item(11,5).Value = 3286  # E11
$ws.Cells.Item(11,6).Value = 'Flattoppad klubbsvamp'  # F11
$ws.Cells.Item(11,7).Value = 'Clavariadelphus truncatus'  # G11
$ws.Cells.Item(11,8).Value = '(Quél.) Donk'  # H11
$ws.Cells.Item(11,16).Value = 'Kalkberget (Kalkberget), Gstr'  # P11
$ws.Cells.Item(11,17).Value = 575021  # Q11
$ws.Cells.Item(11,18).Value = 6703397  # R11
$ws.Cells.Item(11,49).Value = 'Patric Engfeldt'  # AW11
$ws.Cells.Item(11,50).Value = 'Patric Engfeldt'  # AX11

# Row 12
$ws.Cells.Item(12,1).Value = 111896635  # A12
$ws.Cells.Item(12,2).Value = 90480  # B12
$ws.Cells.Item(12,5).Value = 4769  # E12
$ws.Cells.Item(12,6).Value = 'Svavelriska'  # F12
$ws.Cells.Item(12,7).Value = 'Lactarius scrobiculatus'  # G12
$ws.Cells.Item(12,8).Value = '(Scop.:Fr.) Fr.'  # H12
$ws.Cells.Item(12,17).Value = 575037  # Q12
$ws.Cells.Item(12,18).Value = 6703389  # R12

# Row 13
$ws.Cells.Item(13,1).Value = 111896640  # A13
$ws.Cells.Item(13,2).Value = 90480  # B13
$ws.Cells.Item(13,5).Value = 4769  # E13
$ws.Cells.Item(13,6).Value = 'Svavelriska'  # F13
$ws.Cells.Item(13,7).Value = 'Lactarius scrobiculatus'  # G13
$ws.Cells.Item(13,8).Value = '(Scop.:Fr.) Fr.'  # H13
$ws.Cells.Item(13,17).Value = 575025  # Q13
$ws.Cells.Item(13,18).Value = 6703369  # R13

# Row 14
$ws.Cells.Item(14,1).Value = 111896641  # A14
$ws.Cells.Item(14,2).Value = 90480  # B14
$ws.Cells.Item(14,17).Value = 575021  # Q14
$ws.Cells.Item(14,18).Value = 6703371  # R14

# Row 15
$ws.Cells.Item(15,1).Value = 111896638  # A15
$ws.Cells.Item(15,2).Value = 90480  # B15
$ws.Cells.Item(15,17).Value = 575087  # Q15
$ws.Cells.Item(15,18).Value = 6703393  # R15

# Row 16
$ws.Cells.Item(16,1).Value = 111884093  # A16
$ws.Cells.Item(16,2).Value = 98980  # B16
$ws.Cells.Item(16,5).Value = 222498  # E16
$ws.Cells.Item(16,6).Value = 'Blåsippa'  # F16
$ws.Cells.Item(16,7).Value = 'Hepatica nobilis'  # G16
$ws.Cells.Item(16,8).Value = 'Schreb.'  # H16
$ws.Cells.Item(16,16).Value = 'Kopparåsen (Kopparåsen), Gstr'  # P16
$ws.Cells.Item(16,17).Value = 575066  # Q16
$ws.Cells.Item(16,18).Value = 6703388  # R16
$ws.Cells.Item(16,49).Value = 'Patric Engfeldt'  # AW16
$ws.Cells.Item(16,50).Value = 'Patric Engfeldt'  # AX16

# Row 17
$ws.Cells.Item(17,1).Value = 111884133  # A17
$ws.Cells.Item(17,2).Value = 89047  # B17
$ws.Cells.Item(17,4).Value = 'NT'  # D17
$ws.Cells.Item(17,5).Value = 3286  # E17
$ws.Cells.Item(17,6).Value = 'Flattoppad klubbsvamp'  # F17
$ws.Cells.Item(17,7).Value = 'Clavariadelphus truncatus'  # G17
$ws.Cells.Item(17,8).Value = '(Quél.) Donk'  # H17
$ws.Cells.Item(17,16).Value = 'Kalkberget (Kalkberget), Gstr'  # P17
$ws.Cells.Item(17,17).Value = 575059  # Q17
$ws.Cells.Item(17,49).Value = 'Patric Engfeldt'  # AW17
$ws.Cells.Item(17,50).Value = 'Patric Engfeldt'  # AX17

# Row 18
$ws.Cells.Item(18,1).Value = 111896639  # A18
$ws.Cells.Item(18,2).Value = 90480  # B18
$ws.Cells.Item(18,16).Value = 'Kratte masugn, Gstr'  # P18
$ws.Cells.Item(18,17).Value = 575089  # Q18
$ws.Cells.Item(18,18).Value = 6703380  # R18
$ws.Cells.Item(18,49).Value = 'Philipp Weiss'  # AW18
$ws.Cells.Item(18,50).Value = 'Philipp Weiss'  # AX18

# Row 19
$ws.Cells.Item(19,1).Value = 111896655  # A19
$ws.Cells.Item(19,2).Value = 89331  # B19
$ws.Cells.Item(19,5).Value = 3215  # E19
$ws.Cells.Item(19,6).Value = 'Rödgul trumpetsvamp'  # F19
$ws.Cells.Item(19,7).Value = 'Craterellus lutescens'  # G19
$ws.Cells.Item(19,8).Value = '(Fr.) Fr.'  # H19
$ws.Cells.Item(19,17).Value = 575105  # Q19
$ws.Cells.Item(19,18).Value = 6703429  # R19

# Row 20
$ws.Cells.Item(20,1).Value = 111896654  # A20
$ws.Cells.Item(20,2).Value = 89331  # B20
$ws.Cells.Item(20,5).Value = 3215  # E20
$ws.Cells.Item(20,6).Value = 'Rödgul trumpetsvamp'  # F20
$ws.Cells.Item(20,7).Value = 'Craterellus lutescens'  # G20
$ws.Cells.Item(20,8).Value = '(Fr.) Fr.'  # H20
$ws.Cells.Item(20,17).Value = 575073  # Q20
$ws.Cells.Item(20,18).Value = 6703422  # R20

# Row 21
$ws.Cells.Item(21,1).Value = 111896643  # A21
$ws.Cells.Item(21,2).Value = 90480  # B21
$ws.Cells.Item(21,4).Value = 'LC'  # D21
$ws.Cells.Item(21,5).Value = 4769  # E21
$ws.Cells.Item(21,6).Value = 'Svavelriska'  # F21
$ws.Cells.Item(21,7).Value = 'Lactarius scrobiculatus'  # G21
$ws.Cells.Item(21,8).Value = '(Scop.:Fr.) Fr.'  # H21
$ws.Cells.Item(21,16).Value = 'Kratte masugn, Gstr'  # P21
$ws.Cells.Item(21,17).Value = 575039  # Q21
$ws.Cells.Item(21,18).Value = 6703416  # R21
$ws.Cells.Item(21,49).Value = 'Philipp Weiss'  # AW21
$ws.Cells.Item(21,50).Value = 'Philipp Weiss'  # AX21

# Row 22
$ws.Cells.Item(22,1).Value = 111896636  # A22
$ws.Cells.Item(22,2).Value = 90480  # B22
$ws.Cells.Item(22,5).Value = 4769  # E22
$ws.Cells.Item(22,6).Value = 'Svavelriska'  # F22
$ws.Cells.Item(22,7).Value = 'Lactarius scrobiculatus'  # G22
$ws.Cells.Item(22,8).Value = '(Scop.:Fr.) Fr.'  # H22
$ws.Cells.Item(22,17).Value = 575109  # Q22
$ws.Cells.Item(22,18).Value = 6703418  # R22

# Row 23
$ws.Cells.Item(23,2).Value = 90835  # B23

# Row 24
$ws.Cells.Item(24,2).Value = 89114  # B24
